$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90edc99e5abc970ae912e370a2b6f7341213eb6c/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90edc99e5abc970ae912e370a2b6f7341213eb6c/e2e/b.md"

# ---------------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) both flip to "Handed back"
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns on the Overview sheet.
$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# ---------------------------------------------------------------------------
# zh-cn sheet: handback happened -> Status / Target File / Handback File /
# Handback DateTime all get populated.
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-29 14:42:22"
$zhcn.Range("K3").Value = "2016-08-29 14:42:22"

# Recreate the hyperlinks so the two new "Latest Target File" links
# (I2, I3) are interleaved with the existing A2/A3 ones in document order.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $aMdUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aMdUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $bMdUrl, "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aMdUrl, "", "", "a.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(10).ColumnWidth = 39.1667

# ---------------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but with its own handback file
# name / datetime.
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-08-29 14:42:29"
$dede.Range("K3").Value = "2016-08-29 14:42:29"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $aMdUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $aMdUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $bMdUrl, "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aMdUrl, "", "", "a.md")

$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(10).ColumnWidth = 39.1667
